# Updated cryptos list on Sun Aug 11 20:42:10 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $text)
    # Force Excel to store the value as text (not auto-convert numeric-looking
    # strings like "505.34" or "6.19" into real numbers), while keeping the
    # cell's original style (no number-format / quote-prefix residue).
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "58.527.43"
Set-TextCell "E2" "  -3.82%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.556.05"
Set-TextCell "E3" "  -1.41%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.06%  "

# Row 5 - BNB
Set-TextCell "D5" "505.34"
Set-TextCell "E5" "  -3.26%  "

# Row 6 - Solana
Set-TextCell "D6" "144.94"
Set-TextCell "E6" "  -5.95%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.22%  "

# Row 8 - XRP
Set-TextCell "D8" "0.556"
Set-TextCell "E8" "  -6.38%  "

# Row 9 - LidoStakedEther
Set-TextCell "D9" "2.550.02"
Set-TextCell "E9" "  -1.85%  "

# Row 10 - Toncoin
Set-TextCell "D10" "6.19"

# Row 11 - Dogecoin
Set-TextCell "E11" "  -3.05%  "

# Row 12 - Cardano
Set-TextCell "E12" "  -4.34%  "

# Row 13 - TRON
Set-TextCell "E13" "  -0.94%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell "D14" "3.010.57"
Set-TextCell "E14" "  -1.20%  "

# Row 15 - WrappedBTC
Set-TextCell "D15" "58.551.13"
Set-TextCell "E15" "  -3.81%  "

# Row 16 - Avalanche
Set-TextCell "D16" "20.54"
Set-TextCell "E16" "  -5.12%  "

# Row 17 - ShibaInu
Set-TextCell "E17" "  -4.82%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.562.74"
Set-TextCell "E18" "  -1.42%  "

# Row 19 - Polkadot
Set-TextCell "E19" "  -4.80%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "334.24"
Set-TextCell "E20" "  -5.38%  "

# Row 21 - Chainlink
Set-TextCell "D21" "10.08"
Set-TextCell "E21" "  -4.66%  "

# Row 22 - Dai
Set-TextCell "E22" "  -0.12%  "

# Row 23 - Uniswap
Set-TextCell "E23" "  -4.64%  "

# Row 24 - Litecoin
Set-TextCell "D24" "59.53"
Set-TextCell "E24" "  -2.08%  "

# Row 25 - Polygon
Set-TextCell "E25" "  -4.57%  "

# Row 26 - Binance-PegBSC-USD
Set-TextCell "E26" "  +0.15%  "

# Row 27 - Kaspa
Set-TextCell "E27" "  -6.71%  "

# Row 28 - PEPE
Set-TextCell "D28" "0.0₃0777"
Set-TextCell "E28" "  -7.90%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell "D29" "6.86"
Set-TextCell "E29" "  -6.97%  "

# Row 30 - USDe
Set-TextCell "E30" "  +0.08%  "

# Row 31 - Aptos
Set-TextCell "D31" "5.85"
Set-TextCell "E31" "  -7.47%  "

# Row 32 - EthereumClassic
Set-TextCell "E32" "  -3.91%  "

# Row 33 - Monero
Set-TextCell "D33" "148.71"
Set-TextCell "E33" "  +0.46%  "

# Row 34 - PancakeSwap
Set-TextCell "E34" "  -3.91%  "

# Row 35 - NEARProtocol
Set-TextCell "E35" "  -6.97%  "

# Row 36 - SuiNetwork
Set-TextCell "D36" "0.902"
Set-TextCell "E36" "  -3.66%  "

# Row 37 - ImmutableX
Set-TextCell "E37" "  -7.73%  "

# Row 38 - OKB
Set-TextCell "D38" "35.85"
Set-TextCell "E38" "  -1.64%  "

# Row 39 - Fetch.AI
Set-TextCell "D39" "0.817"
Set-TextCell "E39" "  -5.11%  "

# Rows 40 & 41 swap: Filecoin <-> Stacks (with updated values)
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D40" "1.38"
Set-TextCell "E40" "  -8.33%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D41" "3.52"
Set-TextCell "E41" "  -7.24%  "

# Row 42 - Bittensor
Set-TextCell "D42" "280.55"
Set-TextCell "E42" "  -2.74%  "

# Row 43 - FirstDigitalUSD
Set-TextCell "D43" "0.998"
Set-TextCell "E43" "  +0.03%  "

# Row 44 - Stellar
Set-TextCell "D44" "0.0981"
Set-TextCell "E44" "  -3.28%  "

# Row 45 - Mantle
Set-TextCell "D45" "0.605"
Set-TextCell "E45" "  -2.54%  "

# Row 46 - Hedera
Set-TextCell "D46" "0.0531"
Set-TextCell "E46" "  -5.01%  "

# Row 47 - WhiteBITCoin
Set-TextCell "D47" "10.32"
Set-TextCell "E47" "  +0.02%  "

# Row 48 - EnergySwap
Set-TextCell "D48" "18.63"
Set-TextCell "E48" "  -4.83%  "

# Row 49 - VeChain
Set-TextCell "D49" "0.0227"
Set-TextCell "E49" "  -4.74%  "

# Rows 50 & 51 swap: Maker <-> RenderToken (with updated values)
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D50" "4.51"
Set-TextCell "E50" "  -7.53%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D51" "1.912.15"
Set-TextCell "E51" "  -2.36%  "
